$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 180 - Jon (person unchanged)
$ws.Range("D180").Value = 42
$ws.Range("E180").Value = 15
$ws.Range("F180").Value = 57
$ws.Range("G180").Value = 143350
$ws.Range("H180").Value = 110
$ws.Range("I180").Value = 0

# Row 181 - now Maisy
$ws.Range("B181").Value = "Maisy"
$ws.Range("D181").Value = 44
$ws.Range("E181").Value = 7
$ws.Range("F181").Value = 51
$ws.Range("G181").Value = 139500
$ws.Range("H181").Value = 110
$ws.Range("I181").Value = 20
$ws.Range("K181").Value = 360

# Row 182 - now Anthony
$ws.Range("B182").Value = "Anthony"
$ws.Range("D182").Value = 45
$ws.Range("E182").Value = 4
$ws.Range("F182").Value = 49
$ws.Range("G182").Value = 150600
$ws.Range("H182").Value = 190
$ws.Range("I182").Value = 110
$ws.Range("K182").Value = 350

# Row 183 - now Matt
$ws.Range("B183").Value = "Matt"
$ws.Range("D183").Value = 41
$ws.Range("E183").Value = 6
$ws.Range("F183").Value = 47
$ws.Range("G183").Value = 130500
$ws.Range("H183").Value = 70
$ws.Range("I183").Value = -40
$ws.Range("K183").Value = 362

# Row 184 - Pepe (person unchanged)
$ws.Range("D184").Value = 34
$ws.Range("F184").Value = 37
$ws.Range("G184").Value = 104600
$ws.Range("I184").Value = 0

# Row 185 - Richard (person unchanged)
$ws.Range("D185").Value = 31
$ws.Range("E185").Value = 3
$ws.Range("F185").Value = 34
$ws.Range("G185").Value = 102050
$ws.Range("I185").Value = -40

# Row 186 - Andy (person unchanged)
$ws.Range("D186").Value = 24
$ws.Range("E186").Value = 7
$ws.Range("F186").Value = 31
$ws.Range("G186").Value = 90750
$ws.Range("I186").Value = -20

# Row 189 - Mark (person unchanged)
$ws.Range("D189").Value = 22
$ws.Range("F189").Value = 24
$ws.Range("G189").Value = 67400
$ws.Range("I189").Value = -70
